# Update the "想去人数" (want-to-go count) figures in column F on both the
# "展览" sheet and the aggregated "全部类型" sheet, per the latest scrape.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 273
$ws1.Range("F3").Value = 76
$ws1.Range("F5").Value = 7502
$ws1.Range("F6").Value = 5529
$ws1.Range("F8").Value = 69
$ws1.Range("F11").Value = 244
$ws1.Range("F12").Value = 173
$ws1.Range("F13").Value = 47

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 273
$ws4.Range("F3").Value = 76
$ws4.Range("F5").Value = 7502
$ws4.Range("F6").Value = 5529
$ws4.Range("F8").Value = 69
$ws4.Range("F11").Value = 244
$ws4.Range("F14").Value = 173
$ws4.Range("F15").Value = 47
